# Auto-generated edit script: updates live market-price derived columns (H-N)
# across 8 job-leve-profit worksheets, refreshed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 104.29412
$ws.Range("I33").Value = 104.29412
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 104.29412
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 124.70588
$ws.Range("N33").ClearContents()

$ws.Range("H87").Value = 41440.5
$ws.Range("J87").Value = 41440.5
$ws.Range("L87").Value = 41440.5
$ws.Range("N87").Value = -43936.5

$ws.Range("H90").Value = 41440.5
$ws.Range("J90").Value = 41440.5
$ws.Range("L90").Value = 124321.5
$ws.Range("N90").Value = -136801.5

$ws.Range("H111").Value = 1762.875
$ws.Range("I111").Value = 1977.3846
$ws.Range("K111").Value = 5932.1538
$ws.Range("M111").Value = -2865.1538

$ws.Range("H112").Value = 1073.2727
$ws.Range("J112").Value = 1073.2727
$ws.Range("L112").Value = 3219.8181
$ws.Range("N112").Value = -5435.8181

$ws.Range("H129").Value = 3291.6667
$ws.Range("I129").Value = 263.5
$ws.Range("J129").Value = 5310.4443
$ws.Range("K129").Value = 790.5
$ws.Range("L129").Value = 15931.3329
$ws.Range("M129").Value = 4209.5
$ws.Range("N129").Value = -25931.3329

$ws.Range("H131").Value = 3166.6667
$ws.Range("I131").Value = 3250
$ws.Range("K131").Value = 9750
$ws.Range("M131").Value = -4710

$ws.Range("H138").Value = 2128.8374
$ws.Range("I138").Value = 1816.8823
$ws.Range("J138").Value = 2213.0159
$ws.Range("K138").Value = 5450.6469
$ws.Range("L138").Value = 6639.047699999999
$ws.Range("M138").Value = -310.6468999999997
$ws.Range("N138").Value = -16919.0477

$ws.Range("H141").Value = 1191.2778
$ws.Range("I141").Value = 966.2727
$ws.Range("K141").Value = 2898.8181
$ws.Range("M141").Value = 2281.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4811.603
$ws.Range("I32").Value = 4946.981
$ws.Range("J32").Value = 4094.1
$ws.Range("K32").Value = 4946.981
$ws.Range("L32").Value = 4094.1
$ws.Range("M32").Value = -4659.981
$ws.Range("N32").Value = -4668.1

$ws.Range("H61").Value = 1996.1538
$ws.Range("I61").Value = 1787.5
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 1787.5
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -1575.5
$ws.Range("N61").Value = -4924

$ws.Range("H74").Value = 111111920
$ws.Range("I74").Value = 125000800
$ws.Range("K74").Value = 125000800
$ws.Range("M74").Value = -124999926

$ws.Range("H77").Value = 111111920
$ws.Range("I77").Value = 125000800
$ws.Range("K77").Value = 625004000
$ws.Range("M77").Value = -624999632

$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

$ws.Range("H136").Value = 1996.1538
$ws.Range("I136").Value = 1787.5
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 5362.5
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -2812.5
$ws.Range("N136").Value = -18600

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H86").Value = 1911.7307
$ws.Range("I86").Value = 1800.2174
$ws.Range("J86").Value = 2766.6667
$ws.Range("K86").Value = 1800.2174
$ws.Range("L86").Value = 2766.6667
$ws.Range("M86").Value = -677.2174
$ws.Range("N86").Value = -5012.6667

$ws.Range("H89").Value = 1911.7307
$ws.Range("I89").Value = 1800.2174
$ws.Range("J89").Value = 2766.6667
$ws.Range("K89").Value = 9001.087
$ws.Range("L89").Value = 13833.3335
$ws.Range("M89").Value = -3385.087
$ws.Range("N89").Value = -25065.3335

$ws.Range("H128").Value = 2508
$ws.Range("I128").Value = 2508
$ws.Range("K128").Value = 7524
$ws.Range("M128").Value = -5034

$ws.Range("H134").Value = 7260.8
$ws.Range("I134").Value = 7762.4614
$ws.Range("K134").Value = 23287.3842
$ws.Range("M134").Value = -20752.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14443.606
$ws.Range("I31").Value = 21179.4
$ws.Range("J31").Value = 4080.8462
$ws.Range("K31").Value = 21179.4
$ws.Range("L31").Value = 4080.8462
$ws.Range("M31").Value = -20884.4
$ws.Range("N31").Value = -4670.8462

$ws.Range("H34").Value = 14443.606
$ws.Range("I34").Value = 21179.4
$ws.Range("J34").Value = 4080.8462
$ws.Range("K34").Value = 21179.4
$ws.Range("L34").Value = 4080.8462
$ws.Range("M34").Value = -20977.4
$ws.Range("N34").Value = -4484.8462

$ws.Range("H58").Value = 23104.826
$ws.Range("I58").Value = 1468.6875
$ws.Range("J58").Value = 72558.86
$ws.Range("K58").Value = 1468.6875
$ws.Range("L58").Value = 72558.86
$ws.Range("M58").Value = -1265.6875
$ws.Range("N58").Value = -72964.86

$ws.Range("H99").Value = 15628529
$ws.Range("I99").Value = 3319.5715
$ws.Range("J99").Value = 45458470
$ws.Range("K99").Value = 3319.5715
$ws.Range("L99").Value = 45458470
$ws.Range("M99").Value = -1821.5715
$ws.Range("N99").Value = -45461466

$ws.Range("H105").Value = 11364432
$ws.Range("I105").Value = 12500775
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 12500775
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = -12499028
$ws.Range("N105").Value = -4494

$ws.Range("H126").Value = 15628529
$ws.Range("I126").Value = 3319.5715
$ws.Range("J126").Value = 45458470
$ws.Range("K126").Value = 9958.7145
$ws.Range("L126").Value = 136375410
$ws.Range("M126").Value = -7488.7145
$ws.Range("N126").Value = -136380350

$ws.Range("H132").Value = 22167.23
$ws.Range("I132").Value = 23623.305
$ws.Range("J132").Value = 11004
$ws.Range("K132").Value = 70869.91500000001
$ws.Range("L132").Value = 33012
$ws.Range("M132").Value = -68339.91500000001
$ws.Range("N132").Value = -38072

$ws.Range("H136").Value = 23104.826
$ws.Range("I136").Value = 1468.6875
$ws.Range("J136").Value = 72558.86
$ws.Range("K136").Value = 4406.0625
$ws.Range("L136").Value = 217676.58
$ws.Range("M136").Value = -1856.0625
$ws.Range("N136").Value = -222776.58

$ws.Range("H140").Value = 49832.5
$ws.Range("J140").Value = 49832.5
$ws.Range("L140").Value = 49832.5
$ws.Range("N140").Value = -60192.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 678.8077
$ws.Range("J122").Value = 845.7059
$ws.Range("L122").Value = 7611.3531
$ws.Range("N122").Value = -12511.3531

$ws.Range("H131").Value = 818.77
$ws.Range("J131").Value = 818.9596
$ws.Range("L131").Value = 2456.8788
$ws.Range("N131").Value = -12536.8788

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 18277.666
$ws.Range("I132").Value = 3608.9048
$ws.Range("J132").Value = 43948
$ws.Range("K132").Value = 10826.7144
$ws.Range("L132").Value = 131844
$ws.Range("M132").Value = -8296.714399999999
$ws.Range("N132").Value = -136904

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5248.625
$ws.Range("J22").Value = 5400
$ws.Range("L22").Value = 5400
$ws.Range("N22").Value = -5990

$ws.Range("H27").Value = 5248.625
$ws.Range("J27").Value = 5400
$ws.Range("L27").Value = 5400
$ws.Range("N27").Value = -5614

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H100").Value = 2199.7
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 2499.5715
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 2499.5715
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -3581.5715

$ws.Range("H132").Value = 1896.7916
$ws.Range("I132").Value = 1564.2106
$ws.Range("J132").Value = 3160.6
$ws.Range("K132").Value = 4692.6318
$ws.Range("L132").Value = 9481.799999999999
$ws.Range("M132").Value = -2162.6318
$ws.Range("N132").Value = -14541.8

$ws.Range("H136").Value = 35330.332
$ws.Range("I136").Value = 56949.445
$ws.Range("J136").Value = 2901.6667
$ws.Range("K136").Value = 170848.335
$ws.Range("L136").Value = 8705.000100000001
$ws.Range("M136").Value = -168298.335
$ws.Range("N136").Value = -13805.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800

$ws.Range("H132").Value = 873.9729599999999
$ws.Range("I132").Value = 587.1070999999999
$ws.Range("J132").Value = 1766.4445
$ws.Range("K132").Value = 1761.3213
$ws.Range("L132").Value = 5299.333500000001
$ws.Range("M132").Value = 768.6787000000002
$ws.Range("N132").Value = -10359.3335

$ws.Range("H136").Value = 38463396
$ws.Range("I136").Value = 47620812
$ws.Range("J136").Value = 2260
$ws.Range("K136").Value = 142862436
$ws.Range("L136").Value = 6780
$ws.Range("M136").Value = -142859886
$ws.Range("N136").Value = -11880
